$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F3").Value = 1716
    $ws.Range("F4").Value = 789
    $ws.Range("F7").Value = 11989
    $ws.Range("F9").Value = 98
    $ws.Range("F11").Value = 413
    $ws.Range("F13").Value = 863
    $ws.Range("F14").Value = 13479
    $ws.Range("C15").Value = "苏州·I COME ACG动漫品牌博览会x中国国际动漫节cosplay超级盛典江苏赛区"
    $ws.Range("F15").Value = 13467
    $ws.Range("F23").Value = 369
    $ws.Range("F24").Value = 173
}
